$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quick data cleaning: remove the erroneous "RBC / Deposit / 158.99" entry
# (row 18) entirely; Excel shifts every row below it up by one so the
# previously-missing negative expense rows (Dominos, Chatters, Superstore,
# McDonald's, Safeway) now occupy rows 18-22.
$ws.Rows.Item(18).Delete() | Out-Null

# Reflect the final selection state captured in the saved workbook.
$ws.Range("A18:XFD18").Select() | Out-Null
